$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CATEGORIA D – SIN VENTAS")

# The edit cyclically rotates product rows 2 -> 6 -> 5 -> 4 -> 2 (i.e. the
# row that used to be row 4 becomes row 2, row 5 becomes row 4, row 6
# becomes row 5, and the original row 2 becomes row 6), and swaps rows 9
# and 10. Using Range.Copy (rather than re-typing .Value) preserves the
# original cell types/styles exactly (e.g. numeric-looking article codes
# stay text, not auto-converted to numbers).

# Stash originals that would otherwise be overwritten before we can reuse them.
$ws.Range("A2:W2").Copy($ws.Range("A100:W100"))
$ws.Range("A9:W9").Copy($ws.Range("A101:W101"))

# 4-cycle among rows 2, 4, 5, 6
$ws.Range("A4:W4").Copy($ws.Range("A2:W2"))
$ws.Range("A5:W5").Copy($ws.Range("A4:W4"))
$ws.Range("A6:W6").Copy($ws.Range("A5:W5"))
$ws.Range("A100:W100").Copy($ws.Range("A6:W6"))

# 2-cycle (swap) between rows 9 and 10
$ws.Range("A10:W10").Copy($ws.Range("A9:W9"))
$ws.Range("A101:W101").Copy($ws.Range("A10:W10"))

# Clean up the scratch rows used as temporary holding area.
$ws.Range("A100:W101").ClearContents()
